# Weekly update: a new price-report row is inserted right after the
# existing header/metadata block (row 16), pushing all the previously
# recorded observations down by one row and appending one more row
# (the former last row) at the bottom (row 80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 16; rows 16..79 shift down to 17..80.
$ws.Rows("16:16").Insert()

# Populate the freshly inserted row 16 with this week's new observation.
$ws.Cells.Item(16, 1).Value  = 1
$ws.Cells.Item(16, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value  = 44525
$ws.Cells.Item(16, 5).Value  = 15
$ws.Cells.Item(16, 6).Value  = "Fruta"
$ws.Cells.Item(16, 7).Value  = 100102
$ws.Cells.Item(16, 8).Value  = "Cítricos"
$ws.Cells.Item(16, 9).Value  = 100102004
$ws.Cells.Item(16, 10).Value = "Mandarina"
$ws.Cells.Item(16, 11).Value = "Murcott"
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 250
$ws.Cells.Item(16, 14).Value = 12000
$ws.Cells.Item(16, 15).Value = 13000
$ws.Cells.Item(16, 16).Value = 12500
$ws.Cells.Item(16, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(16, 19).Value = 625
$ws.Cells.Item(16, 20).Value = 20
